# Regen save_data to use K instead of Strike#: update column G (K) values
# for rows 2-16 on the active worksheet with the freshly pulled/regenerated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 4
    4  = 3
    5  = 1
    6  = 8
    7  = 8
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 2
    13 = 2
    14 = 3
    15 = 1
    16 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
